$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Time value for row 2 (was "4/23/2025, 7:45:54 PM")
$ws.Range("D2").Value = "2025-05-13 13:02"

# Remove row 3 entirely (ID 524223430 / eve / evvochkaaaaa / 4/23/2025, 7:48:52 PM)
$ws.Rows("3").Delete()
